# Update to agent strategy
# - Corrects the timestamp on the most recent existing run (row 3)
# - Appends a new run row (row 4) with the latest agent strategy results

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the timestamp on row 3 (tiny floating point correction)
$ws.Range("A3").Value = 45791.61911637732

# Append a new results row for the latest strategy run by copying the
# previous row's contents (preserves shared-string references for the
# text columns) then updating the run timestamp.
$ws.Range("B3:S3").Copy()
$ws.Range("B4").PasteSpecial()

$ws.Range("A4").Value = 45791.64495566366
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat
